$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.394933104515076
$ws.Range("B1").Value = 1.314889311790466
$ws.Range("C1").Value = 3.293345212936401
$ws.Range("D1").Value = 2.821429252624512
$ws.Range("E1").Value = 0.8939221501350403
